# Actualización automática 2025-10-29 12:30:07
#
# Updates figures for "GLOBALMATCH S.A.S." (asesor CASTRO ALCIVAR EDA MARIA)
# across the three report sheets, plus the dependent totals/percentages.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (sales by product group) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 28 = CASTRO ALCIVAR EDA MARIA / GLOBALMATCH S.A.S.
$wsGrupo.Range("D28").Value = 356.16      # 240X80 PORCELANATO
$wsGrupo.Range("K28").Value = 2278.38     # PANELES DECORATIVOS
$wsGrupo.Range("M28").Value = 45.36       # PORCELANATO

# Row 60 = "clients reached / 58" counters
$wsGrupo.Range("D60").Value = "7 de 58"
$wsGrupo.Range("M60").Value = "10 de 58"

# --- Sheet "VENTA MENSUAL" (monthly sales) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Row 28 = CASTRO ALCIVAR EDA MARIA / GLOBALMATCH S.A.S., column F = octubre
$wsMensual.Range("F28").Value = 3503.69

# Row 60 = totals row
$wsMensual.Range("F60").Value = 57759.02

# --- Sheet "CUMPLIMIENTO MENSUAL" (monthly compliance) ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 = 240X80 PORCELANATO
$wsCumpl.Range("D3").Value = 8183.14
$wsCumpl.Range("E3").Value = 12204.3374217135
$wsCumpl.Range("F3").Value = 0.4013807020227338

# Row 10 = PANELES DECORATIVOS
$wsCumpl.Range("D10").Value = 8882.879999999999
$wsCumpl.Range("E10").Value = -6166.124115259259
$wsCumpl.Range("F10").Value = 3.269664400063568

# Row 12 = PORCELANATO
$wsCumpl.Range("D12").Value = 27142.66
$wsCumpl.Range("E12").Value = 21481.4
$wsCumpl.Range("F12").Value = 0.5582145958194359

# Row 14 = TOTAL
$wsCumpl.Range("D14").Value = 63363.18
$wsCumpl.Range("E14").Value = 36534.81284188786
$wsCumpl.Range("F14").Value = 0.6342788097883727
